# Reorganize the "family"/"cmdline" example sheets:
#  - the sheet previously named "family" becomes "cmd" and holds the Stata
#    command name ("mestreg")
#  - "cmdline" keeps holding the full command line (unchanged)
#  - a brand-new "family" sheet is added holding the distribution family
#    ("weibull")
#  - a brand-new "frm" sheet is added holding the model form ("hazard")

$wb = $excel.ActiveWorkbook

# Repurpose the old "family" sheet -> "cmd"
$cmdSheet = $wb.Worksheets.Item("family")
$cmdSheet.Range("A1").Value = "mestreg"
$cmdSheet.Name = "cmd"

# Add the new "family" sheet right after "cmdline"
$cmdlineSheet = $wb.Worksheets.Item("cmdline")
$newFamily = $wb.Worksheets.Add($null, $cmdlineSheet)
$newFamily.Name = "family"
$newFamily.Range("A1").Value = "weibull"

# Add the new "frm" sheet right after the new "family" sheet
$newFrm = $wb.Worksheets.Add($null, $newFamily)
$newFrm.Name = "frm"
$newFrm.Range("A1").Value = "hazard"

# Keep the originally active sheet selected
$wb.Worksheets.Item("e(b)").Activate()
